{"js": "// Convert the plain-text URL in the final paragraph into a real hyperlink,\n// then append the new \"recommendation run\" output paragraphs after it.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetText = \"https://www.analyticsvidhya.com/blog/2016/06/quick-guide-build-recommendation-engine-python/\";\n\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === targetText) {\n    targetParagraph = paragraphs.items[i];\n  }\n}\n\nif (!targetParagraph) {\n  throw new Error(\"Could not locate the analyticsvidhya URL paragraph\");\n}\n\n// Turn the plain-text URL into a hyperlink (this also applies the built-in\n// \"Hyperlink\" character style to the run, matching Word's own behavior).\nconst urlRange = targetParagraph.getRange();\nurlRange.hyperlink = targetText;\nawait context.sync();\n\n// New paragraphs of recommender-output text to insert after the URL paragraph.\nconst newParagraphTexts = [\n    \"Setup complete in  658.4465420246124 time\",\n    \"154  ,score: 0.0568964085528\",\n    \"163  ,score: 0.0568964085528\",\n    \"185  ,score: 0.0568964085528\",\n    \"209  ,score: 0.0568964085528\",\n    \"357  ,score: 0.0568964085528\",\n    \"562  ,score: 0.0568964085528\",\n    \"563  ,score: 0.0568964085528\",\n    \"704  ,score: 0.0568964085528\",\n    \"868  ,score: 0.0568964085528\",\n    \"874  ,score: 0.0568964085528\",\n    \"1104  ,score: 0.0568964085528\",\n    \"1122  ,score: 0.0568964085528\",\n    \"1412  ,score: 0.0568964085528\",\n    \"1414  ,score: 0.0568964085528\",\n    \"1415  ,score: 0.0568964085528\",\n    \"1793  ,score: 0.0568964085528\",\n    \"2444  ,score: 0.0568964085528\",\n    \"2526  ,score: 0.0568964085528\",\n    \"4123  ,score: 0.0568964085528\",\n    \"7792  ,score: 0.0568964085528\",\n    \"9646  ,score: 0.0568964085528\",\n    \"9647  ,score: 0.0568964085528\",\n    \"86  ,score: 0.056236896068\",\n    \"917  ,score: 0.0555773835831\",\n    \"954  ,score: 0.0555773835831\",\n    \"961  ,score: 0.0555773835831\",\n    \"1130  ,score: 0.0555773835831\",\n    \"1274  ,score: 0.0555773835831\",\n    \"1406  ,score: 0.0555773835831\",\n    \"1979  ,score: 0.0555773835831\",\n    \"3227  ,score: 0.0555773835831\",\n    \"3500  ,score: 0.0555773835831\",\n    \"3940  ,score: 0.0555773835831\",\n    \"4917  ,score: 0.0555773835831\",\n    \"7232  ,score: 0.0555773835831\",\n    \"7241  ,score: 0.0555773835831\",\n    \"8281  ,score: 0.0555773835831\",\n    \"8282  ,score: 0.0555773835831\",\n    \"8283  ,score: 0.0555773835831\",\n    \"8284  ,score: 0.0555773835831\",\n    \"8285  ,score: 0.0555773835831\",\n    \"8286  ,score: 0.0555773835831\",\n    \"8288  ,score: 0.0555773835831\",\n    \"8289  ,score: 0.0555773835831\",\n    \"8290  ,score: 0.0555773835831\",\n    \"8291  ,score: 0.0555773835831\",\n    \"2179  ,score: 0.0536534669525\",\n    \"314  ,score: 0.0525669885227\",\n    \"330  ,score: 0.0517295503219\",\n    \"468  ,score: 0.0517295503219\",\n    \"Recommendation complete in  35.18560004234314 time\"\n];\n\nlet insertAfter = targetParagraph;\nfor (const text of newParagraphTexts) {\n  insertAfter = insertAfter.insertParagraph(text, Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "# Convert the plain-text URL in the final paragraph into a real hyperlink,\n# then append the new \"recommendation run\" output paragraphs after it.\n\n$d = $word.ActiveDocument\n$url = \"https://www.analyticsvidhya.com/blog/2016/06/quick-guide-build-recommendation-engine-python/\"\n\n# Locate the paragraph whose text is exactly the bare URL (Paragraph.Range.Text\n# includes the trailing paragraph-mark character, so trim before comparing).\n$target = $null\nforeach ($p in $d.Paragraphs) {\n  if ($p.Range.Text.Trim() -eq $url) {\n    $target = $p\n  }\n}\n\nif ($target -eq $null) {\n  throw \"Could not locate the analyticsvidhya URL paragraph\"\n}\n\n# Turn the plain-text URL into a hyperlink in place (applies the built-in\n# \"Hyperlink\" character style to the run, matching Word's own behavior).\n$target.Range.Hyperlink = $url\n\n# New paragraphs of recommender-output text to insert after the URL paragraph.\n$newParagraphTexts = @(\n    \"Setup complete in  658.4465420246124 time\",\n    \"154  ,score: 0.0568964085528\",\n    \"163  ,score: 0.0568964085528\",\n    \"185  ,score: 0.0568964085528\",\n    \"209  ,score: 0.0568964085528\",\n    \"357  ,score: 0.0568964085528\",\n    \"562  ,score: 0.0568964085528\",\n    \"563  ,score: 0.0568964085528\",\n    \"704  ,score: 0.0568964085528\",\n    \"868  ,score: 0.0568964085528\",\n    \"874  ,score: 0.0568964085528\",\n    \"1104  ,score: 0.0568964085528\",\n    \"1122  ,score: 0.0568964085528\",\n    \"1412  ,score: 0.0568964085528\",\n    \"1414  ,score: 0.0568964085528\",\n    \"1415  ,score: 0.0568964085528\",\n    \"1793  ,score: 0.0568964085528\",\n    \"2444  ,score: 0.0568964085528\",\n    \"2526  ,score: 0.0568964085528\",\n    \"4123  ,score: 0.0568964085528\",\n    \"7792  ,score: 0.0568964085528\",\n    \"9646  ,score: 0.0568964085528\",\n    \"9647  ,score: 0.0568964085528\",\n    \"86  ,score: 0.056236896068\",\n    \"917  ,score: 0.0555773835831\",\n    \"954  ,score: 0.0555773835831\",\n    \"961  ,score: 0.0555773835831\",\n    \"1130  ,score: 0.0555773835831\",\n    \"1274  ,score: 0.0555773835831\",\n    \"1406  ,score: 0.0555773835831\",\n    \"1979  ,score: 0.0555773835831\",\n    \"3227  ,score: 0.0555773835831\",\n    \"3500  ,score: 0.0555773835831\",\n    \"3940  ,score: 0.0555773835831\",\n    \"4917  ,score: 0.0555773835831\",\n    \"7232  ,score: 0.0555773835831\",\n    \"7241  ,score: 0.0555773835831\",\n    \"8281  ,score: 0.0555773835831\",\n    \"8282  ,score: 0.0555773835831\",\n    \"8283  ,score: 0.0555773835831\",\n    \"8284  ,score: 0.0555773835831\",\n    \"8285  ,score: 0.0555773835831\",\n    \"8286  ,score: 0.0555773835831\",\n    \"8288  ,score: 0.0555773835831\",\n    \"8289  ,score: 0.0555773835831\",\n    \"8290  ,score: 0.0555773835831\",\n    \"8291  ,score: 0.0555773835831\",\n    \"2179  ,score: 0.0536534669525\",\n    \"314  ,score: 0.0525669885227\",\n    \"330  ,score: 0.0517295503219\",\n    \"468  ,score: 0.0517295503219\",\n    \"Recommendation complete in  35.18560004234314 time\"\n)\n\n$cur = $target\nforeach ($t in $newParagraphTexts) {\n  $cur.Range.InsertParagraphAfter()\n  $cur = $d.Paragraphs.Last\n  $cur.Range.Text = $t\n}\n"}
